# Merged PR 58: Merge in master_dev
#
# 1) "Internal - Data Validation" sheet: insert a new "Scaled Beta"
#    distribution-type row (row 6), pushing the later rows down by one.
# 2) Extend the two workbook-level defined names that cover that table so
#    they still span the whole (now one-row-taller) table.
# 3) "Extent of Contamination" sheet: the default "Step" parameter for the
#    three area-contaminated rows changes from 10000 to 1000.
# 4) Minor view-state touch-ups (selection) on both touched sheets.

$wb = $excel.ActiveWorkbook

# --- 1) Internal - Data Validation: insert "Scaled Beta" row ---------------
$wsVal = $wb.Worksheets.Item("Internal - Data Validation")
$wsVal.Activate()

# Make room for the new row; formats shift down along with the existing data.
$wsVal.Rows.Item(6).Insert()

# Row 4 ("Truncated Normal") has the same 4-parameter column layout
# (A:F populated, G:H blank) that the new row needs, so borrow its
# formatting for the freshly inserted, currently blank row 6.
$wsVal.Range("A4:H4").Copy()
$wsVal.Range("A6:H6").PasteSpecial(-4122)

$wsVal.Range("A6").Value = "Scaled Beta"
$wsVal.Range("B6").Value = 4
$wsVal.Range("C6").Value = "Min (a)"
$wsVal.Range("D6").Value = "Max (b)"
$wsVal.Range("E6").Value = "Alpha"
$wsVal.Range("F6").Value = "Beta"

$wsVal.Range("F15").Select()

# --- 2) Extend the defined names that cover the distribution-type table ----
$wb.Names.Item("Validation_Distribution_Parameter_Count").RefersTo = "='Internal - Data Validation'!`$A`$2:`$B`$12"
$wb.Names.Item("Validation_Distribution_Types").RefersTo = "='Internal - Data Validation'!`$A`$2:`$A`$12"

# --- 3) Extent of Contamination: Step column 10000 -> 1000 -----------------
$wsExt = $wb.Worksheets.Item("Extent of Contamination")
$wsExt.Activate()

$wsExt.Range("O2").Value = 1000
$wsExt.Range("O4").Value = 1000
$wsExt.Range("O6").Value = 1000

$wsExt.Range("F3").Select()
